$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the YearLevel/Birthday header columns (J1/K1)
$ws.Range("J1").Value = "Birthday"
$ws.Range("K1").Value = "YearLevel"

# Update sample data from "Teacher" to "Student"
$ws.Range("C2").Value = "Student"
$ws.Range("D2").Value = "sample.student"
$ws.Range("F2").Value = "sample.student@gmail.com"

# Swap the FIRST / 03/18/2005 sample values (J2/K2) to match new column order
# Force text format so the date-like string is not auto-converted to a date serial,
# then restore the original (default) number format so no extra style is introduced
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "03/18/2005"
$ws.Range("J2").Style = "Normal"
$ws.Range("K2").Value = "FIRST"
